# Auto-generated edit script: updates mass-flow result values
# in the 'Output_flows' and 'Input_flows' worksheets to reflect
# newly-added input files (per commit message 'added new input files').
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Output_flows")
$ws1.Cells.Item(2, 3).Value = [double]"1.924228400277192E-15"
$ws1.Cells.Item(2, 5).Value = [double]"7.860091959506291E-12"
$ws1.Cells.Item(2, 7).Value = [double]"1.388037385308804E-12"
$ws1.Cells.Item(2, 9).Value = [double]"3.993020435635347E-13"
$ws1.Cells.Item(2, 13).Value = [double]"2.830263750140064E-31"
$ws1.Cells.Item(3, 3).Value = [double]"4.071450000026934E-17"
$ws1.Cells.Item(3, 4).Value = [double]"4.023192276894676E-17"
$ws1.Cells.Item(3, 5).Value = [double]"4.444473039693597E-12"
$ws1.Cells.Item(3, 7).Value = [double]"2.936930362133114E-14"
$ws1.Cells.Item(3, 9).Value = [double]"8.448780326926443E-15"
$ws1.Cells.Item(3, 13).Value = [double]"5.988518485603905E-33"
$ws1.Cells.Item(4, 3).Value = [double]"1.884277583969626E-17"
$ws1.Cells.Item(4, 4).Value = [double]"1.861943784966587E-16"
$ws1.Cells.Item(4, 5).Value = [double]"2.162308381617226E-09"
$ws1.Cells.Item(4, 7).Value = [double]"1.359218963025609E-14"
$ws1.Cells.Item(4, 9).Value = [double]"3.910117373860849E-15"
$ws1.Cells.Item(4, 13).Value = [double]"2.771501834367737E-33"
$ws1.Cells.Item(5, 3).Value = [double]"9.400925192258398E-18"
$ws1.Cells.Item(5, 4).Value = [double]"9.289498736054326E-16"
$ws1.Cells.Item(5, 5).Value = [double]"4.62888654980407E-06"
$ws1.Cells.Item(5, 7).Value = [double]"6.781334077319657E-15"
$ws1.Cells.Item(5, 9).Value = [double]"1.950812408815889E-15"
$ws1.Cells.Item(5, 13).Value = [double]"1.382741143701788E-33"
$ws1.Cells.Item(6, 3).Value = [double]"7.71409413765097E-15"
$ws1.Cells.Item(6, 4).Value = [double]"7.622661203657316E-12"
$ws1.Cells.Item(6, 5).Value = [double]"33.08722462290815"
$ws1.Cells.Item(6, 7).Value = [double]"5.564542678669839E-12"
$ws1.Cells.Item(6, 9).Value = [double]"1.600773355679499E-12"
$ws1.Cells.Item(6, 13).Value = [double]"1.134632510351495E-30"
$ws1.Cells.Item(7, 3).Value = [double]"5.738913367099217E-09"
$ws1.Cells.Item(7, 6).Value = [double]"2.49721447922244E-13"
$ws1.Cells.Item(7, 7).Value = [double]"4.139750927402712E-05"
$ws1.Cells.Item(7, 9).Value = [double]"1.190898042553938E-05"
$ws1.Cells.Item(7, 10).Value = [double]"3.674634116665635E-05"
$ws1.Cells.Item(7, 13).Value = [double]"8.441117731011226E-24"
$ws1.Cells.Item(8, 3).Value = [double]"2.825916744246076E-09"
$ws1.Cells.Item(8, 4).Value = [double]"2.793352464957851E-10"
$ws1.Cells.Item(8, 6).Value = [double]"3.08429679034933E-12"
$ws1.Cells.Item(8, 7).Value = [double]"2.038468036444449E-05"
$ws1.Cells.Item(8, 9).Value = [double]"5.864139260990636E-06"
$ws1.Cells.Item(8, 13).Value = [double]"4.156517865031691E-24"
$ws1.Cells.Item(9, 3).Value = [double]"1.883633402603065E-09"
$ws1.Cells.Item(9, 4).Value = [double]"1.861307859331599E-09"
$ws1.Cells.Item(9, 6).Value = [double]"2.161553663107328E-09"
$ws1.Cells.Item(9, 7).Value = [double]"1.358754284394123E-05"
$ws1.Cells.Item(9, 9).Value = [double]"3.908780615001784E-06"
$ws1.Cells.Item(9, 13).Value = [double]"2.770554336050993E-24"
$ws1.Cells.Item(10, 3).Value = [double]"9.400924593938789E-10"
$ws1.Cells.Item(10, 4).Value = [double]"9.28949814792293E-09"
$ws1.Cells.Item(10, 6).Value = [double]"4.628886173410456E-06"
$ws1.Cells.Item(10, 7).Value = [double]"6.781333645723291E-06"
$ws1.Cells.Item(10, 9).Value = [double]"1.950812284656906E-06"
$ws1.Cells.Item(10, 13).Value = [double]"1.382741055697562E-24"
$ws1.Cells.Item(11, 3).Value = [double]"7.714094137810025E-07"
$ws1.Cells.Item(11, 4).Value = [double]"7.622661203817028E-05"
$ws1.Cells.Item(11, 6).Value = [double]"33.0872246229157"
$ws1.Cells.Item(11, 7).Value = [double]"0.005564542678784574"
$ws1.Cells.Item(11, 9).Value = [double]"0.001600773355712505"
$ws1.Cells.Item(11, 13).Value = [double]"1.13463251037489E-21"
$ws1.Cells.Item(12, 3).Value = [double]"1.385229483765106E-16"
$ws1.Cells.Item(12, 5).Value = [double]"2.022803479585955E-09"
$ws1.Cells.Item(12, 9).Value = [double]"1.437264317458851E-13"
$ws1.Cells.Item(12, 10).Value = [double]"9.609770644213841E-11"
$ws1.Cells.Item(12, 13).Value = [double]"1.018736859138146E-31"
$ws1.Cells.Item(13, 3).Value = [double]"1.04698002184134E-17"
$ws1.Cells.Item(13, 4).Value = [double]"7.759278520994797E-17"
$ws1.Cells.Item(13, 5).Value = [double]"5.05286116388852E-10"
$ws1.Cells.Item(13, 9).Value = [double]"1.086308834832752E-14"
$ws1.Cells.Item(13, 10).Value = [double]"1.420554360296129E-11"
$ws1.Cells.Item(13, 13).Value = [double]"7.699786580718623E-33"
$ws1.Cells.Item(14, 3).Value = [double]"4.064166127650053E-19"
$ws1.Cells.Item(14, 4).Value = [double]"1.204798421447185E-17"
$ws1.Cells.Item(14, 5).Value = [double]"1.483290506068056E-09"
$ws1.Cells.Item(14, 9).Value = [double]"4.216832679318602E-16"
$ws1.Cells.Item(14, 10).Value = [double]"3.198658393985607E-12"
$ws1.Cells.Item(14, 13).Value = [double]"2.988902477475663E-34"
$ws1.Cells.Item(15, 3).Value = [double]"4.472169877905693E-13"
$ws1.Cells.Item(15, 4).Value = [double]"1.12688644708303E-10"
$ws1.Cells.Item(15, 5).Value = [double]"1.334537966929881"
$ws1.Cells.Item(15, 9).Value = [double]"4.640162704057854E-10"
$ws1.Cells.Item(15, 13).Value = [double]"3.28895995092928E-28"
$ws1.Cells.Item(16, 3).Value = [double]"9.697095739591395E-13"
$ws1.Cells.Item(16, 4).Value = [double]"2.400330794247101E-09"
$ws1.Cells.Item(16, 5).Value = [double]"39748.1985815899"
$ws1.Cells.Item(16, 9).Value = [double]"1.006135795753847E-09"
$ws1.Cells.Item(16, 13).Value = [double]"7.131517898148053E-28"
$ws1.Cells.Item(17, 3).Value = [double]"6.312232861210354E-09"
$ws1.Cells.Item(17, 6).Value = [double]"3.68778461926414E-10"
$ws1.Cells.Item(17, 9).Value = [double]"2.619738364288883E-05"
$ws1.Cells.Item(17, 10).Value = [double]"0.0175204350142379"
$ws1.Cells.Item(17, 13).Value = [double]"1.856877681147727E-23"
$ws1.Cells.Item(18, 3).Value = [double]"7.086430908470602E-10"
$ws1.Cells.Item(18, 4).Value = [double]"2.100757180370353E-10"
$ws1.Cells.Item(18, 6).Value = [double]"1.368112088052449E-10"
$ws1.Cells.Item(18, 9).Value = [double]"2.941050389773365E-06"
$ws1.Cells.Item(18, 10).Value = [double]"0.003846332529388074"
$ws1.Cells.Item(18, 13).Value = [double]"2.084624519129565E-24"
$ws1.Cells.Item(19, 3).Value = [double]"8.383612653098849E-11"
$ws1.Cells.Item(19, 4).Value = [double]"9.94109484862929E-11"
$ws1.Cells.Item(19, 6).Value = [double]"1.223904102365238E-09"
$ws1.Cells.Item(19, 9).Value = [double]"3.479414049127687E-07"
$ws1.Cells.Item(19, 10).Value = [double]"0.002639304320130874"
$ws1.Cells.Item(19, 13).Value = [double]"2.466218145815036E-25"
$ws1.Cells.Item(20, 3).Value = [double]"0.0001118042502996722"
$ws1.Cells.Item(20, 4).Value = [double]"0.001126886481222446"
$ws1.Cells.Item(20, 6).Value = [double]"1.334537966922962"
$ws1.Cells.Item(20, 9).Value = [double]"0.4640162843175683"
$ws1.Cells.Item(20, 13).Value = [double]"3.288960049536375E-19"
$ws1.Cells.Item(21, 3).Value = [double]"0.0002424273934949869"
$ws1.Cells.Item(21, 4).Value = [double]"0.02400330794299403"
$ws1.Cells.Item(21, 6).Value = [double]"39748.19858159229"
$ws1.Cells.Item(21, 9).Value = [double]"1.006135795775437"
$ws1.Cells.Item(21, 13).Value = [double]"7.131517898301081E-19"

$ws2 = $wb.Worksheets.Item("Input_flows")
$ws2.Cells.Item(2, 3).Value = [double]"4.063066203009719E-12"
$ws2.Cells.Item(3, 3).Value = [double]"2.543972120027284E-13"
$ws2.Cells.Item(4, 3).Value = [double]"9.973943716433702E-15"
$ws2.Cells.Item(5, 3).Value = [double]"4.931050212011886E-15"
$ws2.Cells.Item(6, 3).Value = [double]"7.251135629953755E-12"
$ws2.Cells.Item(7, 3).Value = [double]"3.669315005443532E-05"
$ws2.Cells.Item(8, 3).Value = [double]"1.481514315796482E-05"
$ws2.Cells.Item(9, 3).Value = [double]"9.875548542039827E-06"
$ws2.Cells.Item(10, 3).Value = [double]"4.931044542685385E-06"
$ws2.Cells.Item(11, 3).Value = [double]"0.007242314063495054"
$ws2.Cells.Item(12, 3).Value = [double]"1.033725256335287E-12"
$ws2.Cells.Item(13, 3).Value = [double]"7.467893116300725E-14"
$ws2.Cells.Item(14, 3).Value = [double]"1.095678471155684E-15"
$ws2.Cells.Item(15, 3).Value = [double]"4.640480103084508E-10"
$ws2.Cells.Item(16, 3).Value = [double]"1.024816990789045E-09"
$ws2.Cells.Item(17, 3).Value = [double]"2.679183820413722E-05"
$ws2.Cells.Item(18, 3).Value = [double]"3.015912258236822E-06"
$ws2.Cells.Item(19, 3).Value = [double]"3.565520086231008E-07"
$ws2.Cells.Item(20, 3).Value = [double]"0.464048028311376"
$ws2.Cells.Item(21, 3).Value = [double]"1.024816990811768"
$ws2.Cells.Item(22, 4).Value = [double]"1.442425064066983E-17"
$ws2.Cells.Item(35, 4).Value = [double]"2.616586762609753E-30"
$ws2.Cells.Item(40, 4).Value = [double]"2.616586757563879E-21"
$ws2.Cells.Item(42, 3).Value = [double]"7.37620448702451E-30"
$ws2.Cells.Item(43, 3).Value = [double]"1.257017738107122E-30"
$ws2.Cells.Item(44, 3).Value = [double]"8.223456996366696E-31"
$ws2.Cells.Item(45, 3).Value = [double]"3.982093810448267E-31"
$ws2.Cells.Item(46, 3).Value = [double]"2.510776958263575E-29"
$ws2.Cells.Item(47, 3).Value = [double]"3.941084024364442E-21"
$ws2.Cells.Item(48, 3).Value = [double]"1.237179805631107E-21"
$ws2.Cells.Item(49, 3).Value = [double]"8.223421897381032E-22"
$ws2.Cells.Item(50, 3).Value = [double]"3.982093953218299E-22"
$ws2.Cells.Item(51, 3).Value = [double]"2.510776958397588E-20"
$ws2.Cells.Item(52, 3).Value = [double]"5.822107125566287E-30"
$ws2.Cells.Item(53, 3).Value = [double]"6.283383363577566E-31"
$ws2.Cells.Item(54, 3).Value = [double]"6.18224932155721E-32"
$ws2.Cells.Item(55, 3).Value = [double]"6.954744104992948E-31"
$ws2.Cells.Item(56, 3).Value = [double]"7.588511955703359E-27"
$ws2.Cells.Item(57, 3).Value = [double]"5.225484583571057E-21"
$ws2.Cells.Item(58, 3).Value = [double]"5.867782217294848E-22"
$ws2.Cells.Item(59, 3).Value = [double]"6.149963049736656E-23"
$ws2.Cells.Item(60, 3).Value = [double]"6.95475965471356E-22"
$ws2.Cells.Item(61, 3).Value = [double]"7.588511956164443E-18"
$ws2.Cells.Item(122, 7).Value = [double]"1.947645331179711E-30"
$ws2.Cells.Item(123, 7).Value = [double]"1.217430609136238E-31"
$ws2.Cells.Item(124, 7).Value = [double]"4.817005613173906E-33"
$ws2.Cells.Item(125, 7).Value = [double]"2.381711900599703E-33"
$ws2.Cells.Item(126, 7).Value = [double]"3.486758553161669E-30"
$ws2.Cells.Item(127, 7).Value = [double]"1.76941283647319E-23"
$ws2.Cells.Item(128, 7).Value = [double]"7.157101744232104E-24"
$ws2.Cells.Item(129, 7).Value = [double]"4.770812870097565E-24"
$ws2.Cells.Item(130, 7).Value = [double]"2.382144122874737E-24"
$ws2.Cells.Item(131, 7).Value = [double]"3.483176506395809E-21"
$ws2.Cells.Item(132, 7).Value = [double]"4.961175708117419E-31"
$ws2.Cells.Item(133, 7).Value = [double]"3.584668136534769E-32"
$ws2.Cells.Item(134, 7).Value = [double]"5.291571705787472E-34"
$ws2.Cells.Item(135, 7).Value = [double]"2.276153373908479E-28"
$ws2.Cells.Item(136, 7).Value = [double]"5.024438084204643E-28"
$ws2.Cells.Item(137, 7).Value = [double]"1.313450951479004E-23"
$ws2.Cells.Item(138, 7).Value = [double]"1.47843215843858E-24"
$ws2.Cells.Item(139, 7).Value = [double]"1.747885183373213E-25"
$ws2.Cells.Item(140, 7).Value = [double]"2.276222738841701E-19"
$ws2.Cells.Item(141, 7).Value = [double]"5.024591008387756E-19"

